# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-name suffixes to "_FV2404" / "_FV2410"
# - Turn the data range into an Excel Table (Table1)
# - Freeze the header row (row 1) and select the top-left cell under it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels ("_old" -> "_FV2404", "_new" -> "_FV2410") ---
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Convert the used range A1:U72 into a native Excel table ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U72"), $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the top row and set the selection under the freeze line ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null
